$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Word" value in B11 from "gyro" to "aa:aa:aa"
$ws.Range("B11").Value = "aa:aa:aa"

# Apply a time number format (h:mm:ss) to B11, matching numFmtId 21
$ws.Range("B11").NumberFormat = "h:mm:ss"

# Update the active selection to B11 (instead of C11)
$ws.Range("B11").Select()

$wb.Save()
